$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "edit1"
$ws.Range("B9").Value = "riya-morankar"
$ws.Range("C9").Value = "Merged"

# E9 holds a date-formatted string ("2025-06-18") that must stay literal
# text (matching the other rows in the Date column) rather than being
# auto-converted into a date serial number by Excel's input parsing.
# Force the cell to text format before assigning the value, then reset
# the style back to Normal so no stray number-format styling is left
# behind on the cell.
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2025-06-18"
$ws.Range("E9").Style = "Normal"

$ws.Range("F9").Value = "N/A"
